# Ran code for averaged intensities on spiral schemes
#
# Adds three new averaging schemes (Gaussian-Quadrature moves up in the
# list, and three new "Spiral-90deg-*" schemes are inserted) to the
# alpha1F averaged-intensities table. The previously-last three schemes
# (NoRotation-tilt60deg, Rotation-NoTilt, Rotation-60detTilt) and the
# HexGrid-* schemes shift down to make room, and the three HexGrid rows
# that fall off the bottom are re-appended as new rows 17-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 10-16: new labels/values due to inserted Spiral schemes ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.004263002895115
$ws.Range("D10").Value = 1.015243375825277
$ws.Range("E10").Value = 0.998206045989162
$ws.Range("F10").Value = 1.004263002895115
$ws.Range("G10").Value = 1.006928091757363
$ws.Range("H10").Value = 0.9995347420427434
$ws.Range("I10").Value = 1.001212588345775
$ws.Range("J10").Value = 1.015243375825277
$ws.Range("K10").Value = 1.006724710907219
$ws.Range("L10").Value = 1.005493856901167
$ws.Range("M10").Value = 1.004231307809239

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.0020827473242
$ws.Range("D11").Value = 1.001629628650072
$ws.Range("E11").Value = 1.002900159458971
$ws.Range("F11").Value = 1.0020827473242
$ws.Range("G11").Value = 1.000792327232479
$ws.Range("H11").Value = 1.007009296195424
$ws.Range("I11").Value = 1.001519202971311
$ws.Range("J11").Value = 1.001629628650072
$ws.Range("K11").Value = 1.002264894054522
$ws.Range("L11").Value = 1.002173820689361
$ws.Range("M11").Value = 1.00265556030541

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.002032971723338
$ws.Range("D12").Value = 1.001692589617002
$ws.Range("E12").Value = 1.002907880122454
$ws.Range("F12").Value = 1.002032971723338
$ws.Range("G12").Value = 1.000773911645138
$ws.Range("H12").Value = 1.00701232816197
$ws.Range("I12").Value = 1.001494237774671
$ws.Range("J12").Value = 1.001692589617002
$ws.Range("K12").Value = 1.002300234869728
$ws.Range("L12").Value = 1.002166603296533
$ws.Range("M12").Value = 1.002652319840762

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.002071122651546
$ws.Range("D13").Value = 1.001621384540416
$ws.Range("E13").Value = 1.002903002257248
$ws.Range("F13").Value = 1.002071122651546
$ws.Range("G13").Value = 1.00076786329088
$ws.Range("H13").Value = 1.006993768153334
$ws.Range("I13").Value = 1.00151139352175
$ws.Range("J13").Value = 1.001621384540416
$ws.Range("K13").Value = 1.002262193398832
$ws.Range("L13").Value = 1.002166658025189
$ws.Range("M13").Value = 1.002644755735862

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9931240000000006
$ws.Range("D14").Value = 1.031467999999999
$ws.Range("E14").Value = 1.002080000000002
$ws.Range("F14").Value = 0.9931240000000006
$ws.Range("G14").Value = 1.018547999999999
$ws.Range("H14").Value = 1.001468
$ws.Range("I14").Value = 1.000415999999999
$ws.Range("J14").Value = 1.031467999999999
$ws.Range("K14").Value = 1.016774000000001
$ws.Range("L14").Value = 1.004949000000001
$ws.Range("M14").Value = 1.007850666666666

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.99
$ws.Range("D15").Value = 1.043887499999999
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.99
$ws.Range("G15").Value = 1.03
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1.043887499999999
$ws.Range("K15").Value = 1.021943749999999
$ws.Range("L15").Value = 1.005971875
$ws.Range("M15").Value = 1.010647916666666

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9954397073408002
$ws.Range("D16").Value = 1.027020472729596
$ws.Range("E16").Value = 1.000843203174398
$ws.Range("F16").Value = 0.9954397073408002
$ws.Range("G16").Value = 1.017484850380798
$ws.Range("H16").Value = 1.001577954918399
$ws.Range("I16").Value = 0.9998255251456001
$ws.Range("J16").Value = 1.027020472729596
$ws.Range("K16").Value = 1.013931837951997
$ws.Range("L16").Value = 1.004685772646399
$ws.Range("M16").Value = 1.007031952281599

# --- Append new rows 17-19 (HexGrid rows moved to bottom), copying style from row 16 ---
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 1.003454490005157
$ws.Range("D17").Value = 1.003892631781617
$ws.Range("E17").Value = 1.002309231329607
$ws.Range("F17").Value = 1.003454490005157
$ws.Range("G17").Value = 1.001763654542156
$ws.Range("H17").Value = 1.003672660080131
$ws.Range("I17").Value = 1.001083391245056
$ws.Range("J17").Value = 1.003892631781617
$ws.Range("K17").Value = 1.003100931555612
$ws.Range("L17").Value = 1.003277710780385
$ws.Range("M17").Value = 1.002696009830621

$ws.Range("A16").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.003253477903471
$ws.Range("D18").Value = 1.002633349208918
$ws.Range("E18").Value = 1.002112904197789
$ws.Range("F18").Value = 1.003253477903471
$ws.Range("G18").Value = 1.000670299236557
$ws.Range("H18").Value = 1.005432656495482
$ws.Range("I18").Value = 1.002227748980802
$ws.Range("J18").Value = 1.002633349208918
$ws.Range("K18").Value = 1.002373126703354
$ws.Range("L18").Value = 1.002813302303412
$ws.Range("M18").Value = 1.00272173933717

$ws.Range("A16").Copy($ws.Range("A19"))
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.004546411594985
$ws.Range("D19").Value = 1.001288074648337
$ws.Range("E19").Value = 1.00297289578985
$ws.Range("F19").Value = 1.004546411594985
$ws.Range("G19").Value = 0.9999812316788993
$ws.Range("H19").Value = 1.005880057747088
$ws.Range("I19").Value = 1.001936564345619
$ws.Range("J19").Value = 1.001288074648337
$ws.Range("K19").Value = 1.002130485219094
$ws.Range("L19").Value = 1.003338448407039
$ws.Range("M19").Value = 1.002767539300796
